# Regenerate the experiment task-order sheets: rename each tab to a freshly
# generated order id and rewrite its stim-file list to match the new order.
# Physical sheet positions (and r:id) are unchanged - only the tab names and
# the B-column contents (plus row counts on sheet 2 and sheet 5) move.

$wb = $excel.ActiveWorkbook

# --- Sheet 1 (was GNG_TO-..., now vSAT_TO-...) ---------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "vSAT_TO-16515889152153957"
$ws1.Range("B2").Value = "SAT_stims-1651588915167187.csv"
$ws1.Range("B3").Value = "vSAT_stims-16515889151969259.csv"
$ws1.Range("B4").Value = "SAT_stims-16515889151442163.csv"
$ws1.Range("B5").Value = "vSAT_stims-16515889151814542.csv"

# --- Sheet 2 (was NB_TO-..., now GNG_TO-...) ------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "GNG_TO-1651588915248668"
$ws2.Range("B2").Value = "go_stims-16515889152204773.csv"
$ws2.Range("B3").Value = "GNG_stims-16515889152306733.csv"
$ws2.Range("B4").Value = "go_stims-1651588915232668.csv"
$ws2.Range("B5").Value = "GNG_stims-1651588915246672.csv"
# rows 6:10 no longer exist for this (shorter) task order
$ws2.Range("A6:B10").Clear()

# --- Sheet 3 (RS_TO-..., stays RS_TO-..., new id) -------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "RS_TO-16515889152516692"
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4 (TOL_TO-..., stays TOL_TO-..., new id) -----------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "TOL_TO-16515889153106687"
$ws4.Range("B2").Value = "MM_stims-16515889152766664.csv"
$ws4.Range("B3").Value = "ZM_stims-16515889152546654.csv"
$ws4.Range("B4").Value = "MM_stims-16515889152926674.csv"
$ws4.Range("B5").Value = "ZM_stims-16515889152776656.csv"
$ws4.Range("B6").Value = "MM_stims-16515889153096704.csv"
$ws4.Range("B7").Value = "ZM_stims-16515889152946746.csv"

# --- Sheet 5 (was vSAT_TO-..., now NB_TO-...) -----------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "NB_TO-16515889165133376"
$ws5.Range("B2").Value = "ZB-match_5-16515889155536628.csv"
$ws5.Range("B3").Value = "OB-16515889156744006.csv"
$ws5.Range("B4").Value = "TB-1651588916204964.csv"
$ws5.Range("B5").Value = "TB-1651588916485304.csv"
# this (longer) task order needs 5 more rows, formatted like the existing ones
$ws5.Range("A6").Value = 4
$ws5.Range("A7").Value = 5
$ws5.Range("A8").Value = 6
$ws5.Range("A9").Value = 7
$ws5.Range("A10").Value = 8
$ws5.Range("B6").Value = "ZB-match_3-16515889153676646.csv"
$ws5.Range("B7").Value = "OB-16515889157094023.csv"
$ws5.Range("B8").Value = "ZB-match_3-1651588915640401.csv"
$ws5.Range("B9").Value = "TB-16515889164402757.csv"
$ws5.Range("B10").Value = "OB-16515889161336493.csv"
# copy column-A's bordered/bold style down onto the newly added rows
$ws5.Range("A5").Copy()
$ws5.Range("A6:A10").PasteSpecial(-4122)
